$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price/volume refresh + two row swaps:
# Litecoin/Polygon order, and VeChain replaced by dogwifhat).
# Note: values that look like a bare number (e.g. "0.999", "146.78")
# are prefixed with a leading apostrophe so Excel stores them as text,
# matching the original inline-string "Price" column formatting instead
# of converting them into numeric cells.

$ws.Range("D2").Value = '62.864.09'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D5").Value = '''583.67'
$ws.Range("E5").Value = '  +2.35%  '
$ws.Range("D6").Value = '''146.78'
$ws.Range("E6").Value = '  +4.61%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("E10").Value = '  +2.40%  '
$ws.Range("D11").Value = '''0.394'
$ws.Range("E11").Value = '  +2.27%  '
$ws.Range("D12").Value = '4.036.59'
$ws.Range("E12").Value = '  +2.00%  '
$ws.Range("D13").Value = '''29.45'
$ws.Range("E13").Value = '  +5.94%  '
$ws.Range("D15").Value = '3.446.96'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("D17").Value = '62.843.31'
$ws.Range("E17").Value = '  +2.86%  '
$ws.Range("D18").Value = '''6.25'
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("D19").Value = '''14.31'
$ws.Range("E19").Value = '  +5.85%  '
$ws.Range("D20").Value = '''9.33'
$ws.Range("E20").Value = '  +5.14%  '
$ws.Range("D21").Value = '''394.35'
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = '''0.563'
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '''75.38'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  +4.47%  '
$ws.Range("D26").Value = '3.587.25'
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("E27").Value = '  -1.72%  '
$ws.Range("D28").Value = '''7.74'
$ws.Range("E28").Value = '  +7.67%  '
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").Value = '''8.20'
$ws.Range("E30").Value = '  +3.26%  '
$ws.Range("E31").Value = '  +7.04%  '
$ws.Range("E32").Value = '  +1.17%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '''23.80'
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("D35").Value = '''5.33'
$ws.Range("E35").Value = '  +7.35%  '
$ws.Range("D36").Value = '''7.08'
$ws.Range("E36").Value = '  +2.42%  '
$ws.Range("E37").Value = '  +9.55%  '
$ws.Range("D38").Value = '''168.23'
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("D39").Value = '''30.45'
$ws.Range("E39").Value = '  +17.11%  '
$ws.Range("D40").Value = '3.476.75'
$ws.Range("E40").Value = '  +1.82%  '
$ws.Range("D41").Value = '''0.0768'
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").Value = '''0.792'
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D43").Value = '''42.86'
$ws.Range("E43").Value = '  +1.28%  '
$ws.Range("E45").Value = '  +5.57%  '
$ws.Range("E46").Value = '  +8.28%  '
$ws.Range("D47").Value = '2.516.99'
$ws.Range("E47").Value = '  +3.38%  '
$ws.Range("D48").Value = '''23.65'
$ws.Range("E48").Value = '  +3.93%  '
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '''2.18'
$ws.Range("E51").Value = '  +4.24%  '
